# The deck ships two theme parts:
#   ppt/theme/theme1.xml  - the "Integral" (Red Violet) theme used by the slide master
#   ppt/theme/theme2.xml  - the default "Office Theme" used by the notes master
# The authored change swaps their contents, so the slide master ends up on the
# plain Office color scheme and the notes master ends up on the Red Violet one.
# Only the 12 theme colors differ between the two parts (font + format schemes
# are identical), so recolor the presentation's theme color scheme to the
# Office palette that used to live in theme2.xml.

$p = $ppt.ActivePresentation

function Make-RGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office theme colour scheme, in ThemeColorScheme.Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = Make-RGB($officeColors[$i - 1])
}
